# Append/update: keep only the newest scraped record (2025-11-23 06:25 JST)
# on the "ランサーズ" sheet, dropping the older rows that were scraped
# at 02:02 and replacing row 2 with the newly captured listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember the existing "Hyperlink" cell style used by F2 so we can
# re-apply it after we rebuild the hyperlink below.
$linkStyle = $ws.Range("F2").Style

# Row 2 becomes the single surviving record with refreshed values.
$ws.Range("A2").Value = "2025-11-23 06:25:17"
$ws.Range("B2").Value = "【緊急】海外からWordPress管理画面にログインできない問題の調査と修正"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = "◇管理 ○WordPress"

# The older records (rows 3-9) are gone in the new snapshot.
$ws.Rows("3:9").Delete()

# Rebuild the hyperlink collection: drop every stale hyperlink (rows 3-9
# no longer exist, and row 2's link must point at the new URL) and add a
# single fresh hyperlink for F2, then restore its original cell style.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5439670") | Out-Null
$ws.Range("F2").Style = $linkStyle
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5439670"

# Column D (価格) is narrower in the new layout.
$ws.Columns("D").ColumnWidth = 25.17
